$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-7 from
# serial 45192 (2023-09-23) to serial 45202 (2023-10-03).
$ws.Range("C2:C7").Value = 45202
